$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.125.77"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "1.905.64"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").Value = "'253.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.36%  "

$ws.Range("D6").Value = "'0.694"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("D8").Value = "'41.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.74%  "

$ws.Range("D9").Value = "'0.358"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.02%  "

$ws.Range("D10").Value = "'52.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("D11").Value = "'0.0751"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.82%  "

$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("D13").Value = "'13.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.84%  "

$ws.Range("D14").Value = "2.181.32"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").Value = "'0.734"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.73%  "

$ws.Range("D16").Value = "'5.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.99%  "

$ws.Range("D17").Value = "1.909.60"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").Value = "35.114.89"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").Value = "'73.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.44%  "

$ws.Range("D20").Value = "0.0₃0838"
$ws.Range("E20").Value = "  +3.07%  "

$ws.Range("D21").Value = "'242.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("D22").Value = "'12.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.58%  "

$ws.Range("D23").Value = "'5.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.12%  "

$ws.Range("E24").Value = "  -0.50%  "

$ws.Range("D25").Value = "'2.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.54%  "

$ws.Range("D26").Value = "'2.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").Value = "'167.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").Value = "'8.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("D29").Value = "'18.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.09%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").Value = "4.128.42"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0605"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.59%  "

$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.61%  "

$ws.Range("D34").Value = "'4.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.96%  "

$ws.Range("D35").Value = "'1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.46%  "

$ws.Range("D36").Value = "'4.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.45%  "

$ws.Range("E37").Value = "  -0.46%  "

$ws.Range("E38").Value = "  -5.41%  "

$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("D40").Value = "'103.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.26%  "

$ws.Range("D41").Value = "'17.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.08%  "

$ws.Range("E42").Value = "  +3.78%  "

$ws.Range("E43").Value = "  +1.80%  "

$ws.Range("D44").Value = "'0.0650"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.24%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.308.89"
$ws.Range("E45").Value = "  -2.31%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'2.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("D47").Value = "'12.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.61%  "

$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("E49").Value = "  -1.17%  "

$ws.Range("E50").Value = "  +2.48%  "

$ws.Range("E51").Value = "  +6.32%  "
